$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "size/" prefix from the image filenames stored in column A
# (e.g. "size/size_0.png" -> "size_0.png") across the used range.
# LookAt = 2 (xlPart) so the match doesn't need to be the whole cell content.
$used = $ws.UsedRange
$used.Replace("size/size_", "size_", 2)
